$d = $word.ActiveDocument

# Locate the target paragraph: it is the tab-only paragraph that immediately
# follows the paragraph containing "tổng 8" (the sentence right before the
# insertion point described by the commit "Chi phí kiểm thử").
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "tổng 8") {
        $candidate = $p.Next()
        if ($candidate.Range.Text -eq "`t`r") {
            $target = $candidate
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph (tab-only paragraph after 'tổng 8')."
}

$r = $target.Range

# Replace that single (tab-only) paragraph with two paragraphs:
#   1) the same tab run, followed by the new "Mỗi yêu cầu test tính ..." sentence
#      (still part of the same paragraph);
#   2) a brand-new paragraph with a tab, a Wingdings arrow symbol, and the
#      "chi phí kiểm thử: 2,4 x 18 = 43,2 triệu" sentence.
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>Mỗi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>yêu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>cầu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> test </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>tính</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">00k -&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>mỗi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>chức</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>năng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">00k x 8 = </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>2,4</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>triệu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>tiền</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>kiểm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>thử</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:i/></w:rPr><w:sym w:font="Wingdings" w:char="F0E8"/></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> chi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>phí</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>kiểm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>thử</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>2,</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">4 x 18 = </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>43,2</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>triệu</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$r.InsertXML($xmlFrag)
Write-Output "Inserted 'Chi phi kiem thu' paragraphs successfully."
